$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 128.28572
$ws.Range("I6").Value = 77.75
$ws.Range("J6").Value = 195.66667
$ws.Range("K6").Value = 233.25
$ws.Range("L6").Value = 587.00001
$ws.Range("M6").Value = -121.25
$ws.Range("N6").Value = -811.00001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34996.668
$ws.Range("J44").Value = 34996.668
$ws.Range("L44").Value = 34996.668
$ws.Range("N44").Value = -35972.668
$ws.Range("H55").Value = 25000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H80").Value = 36811.273
$ws.Range("J80").Value = 39991.555
$ws.Range("L80").Value = 39991.555
$ws.Range("N80").Value = -41987.555
$ws.Range("H81").Value = 39999
$ws.Range("J81").Value = 39999
$ws.Range("L81").Value = 39999
$ws.Range("N81").Value = -41995
$ws.Range("H83").Value = 36811.273
$ws.Range("J83").Value = 39991.555
$ws.Range("L83").Value = 119974.665
$ws.Range("N83").Value = -129958.665
$ws.Range("H84").Value = 39999
$ws.Range("J84").Value = 39999
$ws.Range("L84").Value = 119997
$ws.Range("N84").Value = -129981
$ws.Range("H122").Value = 2932.4285
$ws.Range("I122").Value = 2932.4285
$ws.Range("K122").Value = 8797.2855
$ws.Range("M122").Value = -6347.2855
$ws.Range("H132").Value = 3685.7
$ws.Range("I132").Value = 3607.125
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 10821.375
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8291.375
$ws.Range("N132").Value = -17060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1302.3043
$ws.Range("I20").Value = 1168.5625
$ws.Range("J20").Value = 1608
$ws.Range("K20").Value = 1168.5625
$ws.Range("L20").Value = 1608
$ws.Range("M20").Value = -921.5625
$ws.Range("N20").Value = -2102

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 88.47619
$ws.Range("I7").Value = 44.588234
$ws.Range("J7").Value = 275
$ws.Range("K7").Value = 44.588234
$ws.Range("L7").Value = 275
$ws.Range("M7").Value = 68.411766
$ws.Range("N7").Value = -501
$ws.Range("H31").Value = 4636.5
$ws.Range("I31").Value = 4276.7
$ws.Range("K31").Value = 4276.7
$ws.Range("M31").Value = -3981.7
$ws.Range("H34").Value = 4636.5
$ws.Range("I34").Value = 4276.7
$ws.Range("K34").Value = 4276.7
$ws.Range("M34").Value = -4074.7
$ws.Range("H50").Value = 29980.7
$ws.Range("J50").Value = 29980.7
$ws.Range("L50").Value = 29980.7
$ws.Range("N50").Value = -31230.7
$ws.Range("H58").Value = 3338.8
$ws.Range("J58").Value = 5996.3335
$ws.Range("L58").Value = 5996.3335
$ws.Range("N58").Value = -6402.3335
$ws.Range("H60").Value = 21360.5
$ws.Range("J60").Value = 24984
$ws.Range("L60").Value = 24984
$ws.Range("N60").Value = -26006
$ws.Range("H62").Value = 8846.25
$ws.Range("J62").Value = 7890
$ws.Range("L62").Value = 7890
$ws.Range("N62").Value = -9138
$ws.Range("H65").Value = 8846.25
$ws.Range("J65").Value = 7890
$ws.Range("L65").Value = 39450
$ws.Range("N65").Value = -45690
$ws.Range("H74").Value = 38188.625
$ws.Range("J74").Value = 38188.625
$ws.Range("L74").Value = 38188.625
$ws.Range("N74").Value = -39936.625
$ws.Range("H77").Value = 38188.625
$ws.Range("J77").Value = 38188.625
$ws.Range("L77").Value = 114565.875
$ws.Range("N77").Value = -123301.875
$ws.Range("H132").Value = 2496.6
$ws.Range("I132").Value = 2218.4443
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6655.3329
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4125.3329
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 3338.8
$ws.Range("J136").Value = 5996.3335
$ws.Range("L136").Value = 17989.0005
$ws.Range("N136").Value = -23089.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 262656.88
$ws.Range("J4").Value = 2450
$ws.Range("L4").Value = 7350
$ws.Range("N4").Value = -7574
$ws.Range("H6").Value = 12006.2
$ws.Range("I6").Value = 10.666667
$ws.Range("K6").Value = 32.000001
$ws.Range("M6").Value = 80.999999
$ws.Range("H36").Value = 99.85714
$ws.Range("I36").Value = 99.85714
$ws.Range("K36").Value = 299.57142
$ws.Range("M36").Value = -130.57142
$ws.Range("H37").Value = 97118.57000000001
$ws.Range("J37").Value = 97118.57000000001
$ws.Range("L37").Value = 291355.71
$ws.Range("N37").Value = -291579.71
$ws.Range("H50").Value = 591.6667
$ws.Range("I50").Value = 591.6667
$ws.Range("K50").Value = 1775.0001
$ws.Range("M50").Value = -1294.0001
$ws.Range("H53").Value = 591.6667
$ws.Range("I53").Value = 591.6667
$ws.Range("K53").Value = 1775.0001
$ws.Range("M53").Value = -1294.0001
$ws.Range("H117").Value = 1469.1428
$ws.Range("I117").Value = 862.6667
$ws.Range("J117").Value = 1924
$ws.Range("K117").Value = 2588.0001
$ws.Range("L117").Value = 5772
$ws.Range("M117").Value = 853.9998999999998
$ws.Range("N117").Value = -12656
$ws.Range("H129").Value = 2149.2
$ws.Range("I129").Value = 817.8
$ws.Range("K129").Value = 2453.4
$ws.Range("M129").Value = 2546.6
$ws.Range("H131").Value = 2188.875
$ws.Range("J131").Value = 3074.875
$ws.Range("L131").Value = 9224.625
$ws.Range("N131").Value = -19304.625
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20434.273
$ws.Range("I43").Value = 9500
$ws.Range("J43").Value = 22864.111
$ws.Range("K43").Value = 9500
$ws.Range("L43").Value = 22864.111
$ws.Range("M43").Value = -9349
$ws.Range("N43").Value = -23166.111
$ws.Range("H46").Value = 14248.5
$ws.Range("J46").Value = 19999.8
$ws.Range("L46").Value = 19999.8
$ws.Range("N46").Value = -20311.8
$ws.Range("H107").Value = 150.8
$ws.Range("I107").Value = 150.8
$ws.Range("K107").Value = 150.8
$ws.Range("M107").Value = 1769.2
$ws.Range("H132").Value = 1754
$ws.Range("I132").Value = 1754
$ws.Range("K132").Value = 5262
$ws.Range("M132").Value = -2732

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 349995
$ws.Range("J54").Value = 349995
$ws.Range("L54").Value = 349995
$ws.Range("N54").Value = -351035
$ws.Range("H62").Value = 4333.3335
$ws.Range("H65").Value = 4333.3335
$ws.Range("H126").Value = 1868
$ws.Range("I126").Value = 1868
$ws.Range("K126").Value = 5604
$ws.Range("M126").Value = -3134
$ws.Range("H132").Value = 1104.2858
$ws.Range("J132").Value = 1088.3334
$ws.Range("L132").Value = 3265.0002
$ws.Range("N132").Value = -8325.0002
$ws.Range("H136").Value = 2011.6316
$ws.Range("I136").Value = 2011.6316
$ws.Range("K136").Value = 6034.8948
$ws.Range("M136").Value = -3484.8948
